# "Update countries & provincias Spain"
#
# The source COVID-19 stats were refreshed (new "Casos totales" etc. for a
# handful of countries) and the table re-sorted descending by "Casos
# totales", which shuffled a few rows:
#   - Kazajistan overtakes Kuwait                    (rows 74/75)
#   - Paraguay overtakes Camboya                      (rows 124/125)
#   - El Salvador jumps ahead of Monaco..Guatemala     (rows 131-136)
#   - Guinea Ecuatorial overtakes Namibia (tie, same stats) (rows 167/168)
#   - Fiyi jumps ahead of Antigua y Barbuda..Mongolia   (rows 170-173)
# plus a straightforward stat refresh for Estados Unidos (row 4) and
# Honduras (row 98) that doesn't change their rank, and the "last updated"
# footer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $pais, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Range("A$r").Value = $pais
    $ws.Range("B$r").Value = $b
    $ws.Range("C$r").Value = $c
    $ws.Range("D$r").Value = $d
    $ws.Range("E$r").Value = $e
    $ws.Range("F$r").Value = $f
    $ws.Range("G$r").Value = $g
    $ws.Range("H$r").Value = $h
}

# row 4 — Estados Unidos: updated stats, rank unchanged
Set-Row 4 "Estados Unidos" 367385 381 19671 336838 8879 5 10876

# rows 74/75 — Kazajistan updates and overtakes Kuwait
Set-Row 74 "Kazajistan" 670 8 46 618 16 0 6
Set-Row 75 "Kuwait"     665 0 103 561 20 0 1

# row 98 — Honduras: updated stats, rank unchanged
Set-Row 98 "Honduras" 305 7 6 277 10 0 22

# rows 124/125 — Paraguay updates and overtakes Camboya
Set-Row 124 "Paraguay" 115 2 15 95 1 0 5
Set-Row 125 "Camboya"  114 0 53 61 1 0 0

# rows 131-136 — El Salvador updates and jumps ahead of Monaco..Guatemala,
# which each shift down one place (values unchanged, just relocated)
Set-Row 131 "El Salvador"       78 9 5 69 4 0 4
Set-Row 132 "Monaco"            77 0 4 72 4 0 1
Set-Row 133 "Liechtenstein"     77 0 55 21 0 0 1
Set-Row 134 "Guayana Francesa"  72 0 34 38 1 0 0
Set-Row 135 "Aruba"             71 0 2 69 0 0 0
Set-Row 136 "Guatemala"         70 0 15 52 3 0 3

# rows 167/168 — Guinea Ecuatorial overtakes Namibia (tied stats, pure swap)
Set-Row 167 "Guinea Ecuatorial" 16 0 3 13 0 0 0
Set-Row 168 "Namibia"           16 0 3 13 0 0 0

# rows 170-173 — Fiyi updates and jumps ahead of Antigua y Barbuda..Mongolia,
# which each shift down one place (values unchanged, just relocated)
Set-Row 170 "Fiyi"               15 1 0 15 0 0 0
Set-Row 171 "Antigua y Barbuda"  15 0 0 15 1 0 0
Set-Row 172 "Dominica"           15 0 1 14 0 0 0
Set-Row 173 "Mongolia"           15 0 2 13 0 0 0

# footer timestamp
$ws.Range("A1").Value = "Datos actualizados a 7 de Abril de 2020 a las 05:22"
